$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-parsed as a number by Excel;
# force Text format, assign, then restore the default "Normal" style so the
# resulting cell carries no explicit style index (matches original formatting).
$textCells = @{
    'D4' = '1.002'
    'D5' = '262.85'
    'D6' = '1.001'
    'D7' = '0.5230'
    'D8' = '0.3242'
    'D9' = '0.06772'
    'D10' = '18.71'
    'D11' = '0.7716'
    'D12' = '0.07753'
    'D14' = '88.29'
    'D15' = '5.013'
    'D16' = '1.002'
    'D19' = '0.000007926'
    'D22' = '4.623'
    'D23' = '9.468'
    'D24' = '5.969'
    'D25' = '143.42'
    'D26' = '2.182'
    'D27' = '1.679'
    'D28' = '17.01'
    'D29' = '111.40'
    'D30' = '4.172'
    'D31' = '0.08752'
    'D32' = '4.109'
    'D33' = '0.04817'
    'D34' = '1.129'
    'D35' = '2.874'
    'D36' = '0.7127'
    'D37' = '3.100'
    'D38' = '0.01786'
    'D39' = '2.185'
    'D40' = '0.4848'
    'D41' = '112.48'
    'D42' = '0.8976'
    'D43' = '6.052'
    'D44' = '1.001'
    'D45' = '7.628'
    'D46' = '0.05907'
    'D47' = '0.4148'
    'D48' = '9.043'
    'D49' = '34.90'
    'D50' = '0.1229'
    'D51' = '0.8839'
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Plain text / non-numeric-looking updates (safe to assign directly).
$ws.Range("D2").Value = '26.521.23'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.847.94'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +2.56%  '
$ws.Range("E8").Value = '  +1.31%  '
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").Value = '1.867.46'
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = '26.585.10'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '2.091.38'
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("E26").Value = '  -6.46%  '
$ws.Range("E27").Value = '  +1.73%  '
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  -1.67%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("E51").Value = '  +3.74%  '
